# CEDEN_Benthic_Data_Dictionary.xlsx edit:
# Insert a new row (row 64) into the "CEDEN_Benthic_Data_Dictionary" sheet for the
# "BenthicLabEffort_AgencyCode" field, pushing the existing "DataQuality" /
# "DataQualityIndicator" rows (and the blank rows below them) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CEDEN_Benthic_Data_Dictionary")

# Insert a new blank row above the current row 64 ("DataQuality"); this shifts
# rows 64:73 down to 65:74 and automatically extends the B2:B65 data validation
# range that covered row 64 to B2:B66.
$ws.Rows.Item(64).Insert()

# Copy the formatting of the row above (row 63, which has the same "no match
# found" look the new row needs) onto the newly inserted row.
$ws.Range("A63:F63").Copy()
$ws.Range("A64:F64").PasteSpecial(-4122)

# Fill in the new row's contents.
$ws.Range("A64").Value = "BenthicLabEffort_AgencyCode"
$ws.Range("B64").Value = "text"
$ws.Range("F64").Value = "text"
$ws.Range("C64").FormulaArray = "=IFERROR(INDEX(Data_Dictionary_FromPDF!B:B,E64),""-"")"
$ws.Range("D64").FormulaArray = "=IFERROR(INDEX(Data_Dictionary_FromPDF!D:D,E64),""-"")"
$ws.Range("E64").Formula = "=MATCH(A64,Data_Dictionary_FromPDF!A:A,0)"

[void]$ws.Range("F64").Select()
